$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new row (becomes row 188), pushing old rows 188..286 down to 189..287 ---
$ws.Rows(188).Insert()

$ws.Range("A188").Value2 = 1
$ws.Range("B188").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C188").Value2 = "Arica y Parinacota"
$ws.Range("D188").Value2 = 44806
$ws.Range("E188").Value2 = 15
$ws.Range("F188").Value2 = "Fruta"
$ws.Range("G188").Value2 = 100108
$ws.Range("H188").Value2 = "Tropicales y subtropicales"
$ws.Range("I188").Value2 = 100108006
$ws.Range("J188").Value2 = "Plátano"
$ws.Range("K188").Value2 = "Sin especificar"
$ws.Range("L188").Value2 = "Verde"
$ws.Range("M188").Value2 = 120
$ws.Range("N188").Value2 = 20000
$ws.Range("O188").Value2 = 21000
$ws.Range("P188").Value2 = 20500
$ws.Range("Q188").Value2 = "`$/caja 20 kilos"
$ws.Range("R188").Value2 = "Ecuador"
$ws.Range("S188").Value2 = 1025
$ws.Range("T188").Value2 = 20

# --- Insert second new row (becomes row 234), pushing what is now row 234..287 down to 235..288 ---
$ws.Rows(234).Insert()

$ws.Range("A234").Value2 = 1
$ws.Range("B234").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C234").Value2 = "Arica y Parinacota"
$ws.Range("D234").Value2 = 44761
$ws.Range("E234").Value2 = 15
$ws.Range("F234").Value2 = "Fruta"
$ws.Range("G234").Value2 = 100108
$ws.Range("H234").Value2 = "Tropicales y subtropicales"
$ws.Range("I234").Value2 = 100108006
$ws.Range("J234").Value2 = "Plátano"
$ws.Range("K234").Value2 = "Sin especificar"
$ws.Range("L234").Value2 = "Pintón"
$ws.Range("M234").Value2 = 130
$ws.Range("N234").Value2 = 30000
$ws.Range("O234").Value2 = 32000
$ws.Range("P234").Value2 = 31000
$ws.Range("Q234").Value2 = "`$/caja 20 kilos"
$ws.Range("R234").Value2 = "Ecuador"
$ws.Range("S234").Value2 = 1550
$ws.Range("T234").Value2 = 20

$ws.Range("A1").Select()
